$d = $word.ActiveDocument

# 1. Split the title run into three runs with the same overall text.
$d.Content.Find.Execute(
    "Sequence model for preparation and clean-up of operation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sequence model for preparation for and clean-up of operation",
    2)

# 2. Append " (10 cm from each other)" after "Measuring and placement of ports"
$d.Content.Find.Execute(
    "Measuring and placement of ports",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Measuring and placement of ports (10 cm from each other)",
    2)

# 3. Append " ( tool changing)" after "Operation"
$d.Content.Find.Execute(
    "Operation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Operation ( tool changing)",
    2)
